# Added: meeeting date time
# Adds MeetingDay* / MeetingTime* columns to the "Cell Groups" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cell Groups")

# Pick up the same "required header" look (bold, shaded fill, thin border)
# used by the other starred header cells (e.g. A1 "Name*") for the two new
# header cells, before we overwrite their text.
$ws.Range("A1").Copy($ws.Range("H1"))
$ws.Range("A1").Copy($ws.Range("I1"))

# Data rows for the new MeetingDay* column
$ws.Range("H2").Value = "Thursday"
$ws.Range("H3").Value = "Friday"
$ws.Range("H4").Value = "Saturday"

# Header labels
$ws.Range("H1").Value = "MeetingDay*"
$ws.Range("I1").Value = "MeetingTime*"

# Data rows for the new MeetingTime* column (stored as Excel time serials)
$ws.Range("I2").Value = 0.58333333333333337
$ws.Range("I3").Value = 0.625
$ws.Range("I4").Value = 0.58333333333333337

# Last row
$ws.Range("H5").Value = "Thursday"
$ws.Range("I5").Value = 0.70833333333333304

# Time format for the meeting time column (numFmtId 20 => h:mm)
$ws.Range("I2:I5").NumberFormat = "h:mm"

# Column widths for the new columns
$ws.Columns.Item(8).ColumnWidth = 12.21875
$ws.Columns.Item(9).ColumnWidth = 13.21875

# Restore selection like the source workbook
$ws.Range("F14").Select()
